$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update row 14 (e013 Gun Load) with the new text that adds the "Note:" paragraph
#    about Advancing Fire, and grow the row height from 105 to 150.
$b14 = @'
<Bold>e013 Gun Load</Bold> 
<InlineUIContainer><Button Content='r4.43' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Mark the type of round you want loaded in the main gun before any action begins by clicking the highlighted box on the Tank Card in the correct ammo type box.
<LineBreak/><LineBreak/>
Note: You must choose HE gunload if you want to use Advancing Fire when entering a new battle per 
<InlineUIContainer><Button Content='r22.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
<LineBreak/><LineBreak/> If you do not want to load the gun or have finished selecting the gun load, click image below to continue.
<LineBreak/><LineBreak/>
                                                  <InlineUIContainer><Image Name='c17GunLoad'  Height='80' Width='80'></Image></InlineUIContainer>
'@
$b14 = $b14.TrimEnd("`r","`n")

$ws.Range("B14").Value2 = $b14
$ws.Rows.Item(14).RowHeight = 150

# 2. Insert a new row 32 for the new "e029a Advancing Fire Not Allowed" event,
#    pushing all subsequent rows down by one.
$ws.Rows.Item(32).Insert()

$a32 = "e029a"
$b32 = @'
<Bold>e029a Advancing Fire Not Allowed</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r22.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
Advancing fire is only allowed if there is an HE Gun Load per <InlineUIContainer><Button Content='r9.61' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Choose image below to continue.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='c44AdvanceFireDeny' Height='120' Width='120'></Image></InlineUIContainer> 
'@
$b32 = $b32.TrimEnd("`r","`n")

$ws.Range("A32").Value2 = $a32
$ws.Range("B32").Value2 = $b32
$ws.Rows.Item(32).RowHeight = 120

# 3. Update the sheet view so the new row is visible, matching the author's
#    final scroll position and selection.
$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 29
